$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 239, pushing existing rows 239:253 down to 240:254
$ws.Rows(239).Insert()

# Populate the new row 239 with the new data record
$ws.Cells.Item(239, 1).Value = 4
$ws.Cells.Item(239, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(239, 3).Value = "Los Lagos"
$ws.Cells.Item(239, 4).Value = 44714
$ws.Cells.Item(239, 5).Value = 10
$ws.Cells.Item(239, 6).Value = 100112044
$ws.Cells.Item(239, 7).Value = "Perejil"
$ws.Cells.Item(239, 8).Value = "Sin especificar"
$ws.Cells.Item(239, 9).Value = "Primera"
$ws.Cells.Item(239, 10).Value = 40
$ws.Cells.Item(239, 11).Value = 6000
$ws.Cells.Item(239, 12).Value = 6000
$ws.Cells.Item(239, 13).Value = 6000
$ws.Cells.Item(239, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(239, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(239, 16).Value = 3000
$ws.Cells.Item(239, 17).Value = 2
$ws.Cells.Item(239, 18).Value = "Hortaliza"
